$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: repurpose "dep_peroni_btl_12p_330_vol" row into the new
#     "mod_vol_pna00_impulse_glass_330ml_1_12pack" taxonomy row. The old
#     A2 formula (=C2&E2&G2) is replaced by a literal, F2/G2 are cleared.
$ws.Range("A2").Value = "mod_vol_pna00_impulse_glass_330ml_1_12pack"
$ws.Range("C2").Value = "mod_vol_"
$ws.Range("D2").Value = "PNA00_IMPULSE_GLASS_330ML_1-12PACK"
$ws.Range("E2").Value = "pna00_impulse_glass_330ml_1_12pack"
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()

# --- Row 3: repurpose "dep_peroni_btl_sin_620_vol" row into the new
#     "mod_dist_pna00_impulse_glass_330ml_1_12pack" taxonomy row.
$ws.Range("A3").Value = "mod_dist_pna00_impulse_glass_330ml_1_12pack"
$ws.Range("B3").Value = "distribution"
$ws.Range("C3").Value = "mod_dist_"
$ws.Range("D3").Value = "PNA00_IMPULSE_GLASS_330ML_1-12PACK"
$ws.Range("E3").Value = "pna00_impulse_glass_330ml_1_12pack"
$ws.Range("F3").ClearContents()

# --- Column widths: column A/B split out of the merged 1-2 group, and
#     D/E get wider to fit the new longer taxonomy text.
$ws.Columns.Item(1).ColumnWidth = 42.5
$ws.Columns.Item(4).ColumnWidth = 29
$ws.Columns.Item(5).ColumnWidth = 31.333333333333332

# --- Selection / scroll position: was row-26 selected with the view
#     scrolled to row 10; now a single cell (B4) is selected and the view
#     resets to the top.
$ws.Range("B4").Select()
